# Update instrument/glider reference designators from GP05MOAS-GL002 to
# GP05MOAS-GL364 across the Moorings and Asset_Cal_Info sheets.

$wb = $excel.ActiveWorkbook

$wsMoorings = $wb.Worksheets.Item("Moorings")
$wsAssetCal = $wb.Worksheets.Item("Asset_Cal_Info")

# --- Moorings sheet ---
$wsMoorings.Range("A2").Value = "GP05MOAS-GL364"

# --- Asset_Cal_Info sheet ---
$wsAssetCal.Range("A3").Value = "GP05MOAS-GL364-00-ENG000000"
$wsAssetCal.Range("A4").Value = "GP05MOAS-GL364-01-FLORDM000"
$wsAssetCal.Range("A5").Value = "GP05MOAS-GL364-01-FLORDM000"
$wsAssetCal.Range("A6").Value = "GP05MOAS-GL364-01-FLORDM000"
$wsAssetCal.Range("A7").Value = "GP05MOAS-GL364-01-FLORDM000"
$wsAssetCal.Range("A8").Value = "GP05MOAS-GL364-02-DOSTAM000"
$wsAssetCal.Range("A9").Value = "GP05MOAS-GL364-04-CTDGVM000"

$wb.Save()
